$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.262296333333333
$ws.Range("H2").Value = 9.786889
$ws.Range("I2").Value = 0.01915820289899999
$ws.Range("J2").Value = 0.01915820289899999
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 4.472204835686778
$ws.Range("R2").Value = 40.249843521181
$ws.Range("S2").Value = 0.0002113271069076412
$ws.Range("T2").Value = 0.0002113271069076411
$ws.Range("G3").Value = 3.262296333333333
$ws.Range("H3").Value = 9.786889
$ws.Range("I3").Value = 0.01915820289899999
$ws.Range("J3").Value = 0.01915820289899999
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 303.1121022089774
$ws.Range("R3").Value = 2728.008919880797
$ws.Range("S3").Value = 0.01432309251968322
$ws.Range("T3").Value = 0.01432309251968322
$ws.Range("G4").Value = 3.262296333333333
$ws.Range("H4").Value = 9.786889
$ws.Range("I4").Value = 0.01915820289899999
$ws.Range("J4").Value = 0.01915820289899999
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 96.95050138542534
$ws.Range("R4").Value = 872.554512468828
$ws.Range("S4").Value = 0.004581245654836133
$ws.Range("T4").Value = 0.004581245654836132
$ws.Range("G5").Value = 3.262296333333333
$ws.Range("H5").Value = 9.786889
$ws.Range("I5").Value = 0.01915820289899999
$ws.Range("J5").Value = 0.01915820289899999
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.9002013125163334
$ws.Range("R5").Value = 8.101811812647
$ws.Range("S5").Value = 0.00004253761757299387
$ws.Range("T5").Value = 0.00004253761757299387
$ws.Range("I6").Value = 0.8527862647199704
$ws.Range("J6").Value = 0.8527862647199704
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 199.0705953472801
$ws.Range("R6").Value = 1791.635358125521
$ws.Range("S6").Value = 0.009406772393210848
$ws.Range("T6").Value = 0.009406772393210847
$ws.Range("I7").Value = 0.8527862647199704
$ws.Range("J7").Value = 0.8527862647199704
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.6375617083446159
$ws.Range("T7").Value = 0.6375617083446158
$ws.Range("I8").Value = 0.8527862647199704
$ws.Range("J8").Value = 0.8527862647199704
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 4315.543392826306
$ws.Range("R8").Value = 38839.89053543675
$ws.Range("S8").Value = 0.20392431327451
$ws.Range("T8").Value = 0.20392431327451
$ws.Range("I9").Value = 0.8527862647199704
$ws.Range("J9").Value = 0.8527862647199704
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 40.07052847513634
$ws.Range("R9").Value = 360.6347562762271
$ws.Range("S9").Value = 0.001893470707633726
$ws.Range("T9").Value = 0.001893470707633726
$ws.Range("G10").Value = 21.305189
$ws.Range("H10").Value = 63.915567
$ws.Range("I10").Value = 0.1251171236325075
$ws.Range("J10").Value = 0.1251171236325075
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 29.20677937729366
$ws.Range("R10").Value = 262.861014395643
$ws.Range("S10").Value = 0.001380121084490843
$ws.Range("T10").Value = 0.001380121084490843
$ws.Range("G11").Value = 21.305189
$ws.Range("H11").Value = 63.915567
$ws.Range("I11").Value = 0.1251171236325075
$ws.Range("J11").Value = 0.1251171236325075
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 1979.544457615566
$ws.Range("R11").Value = 17815.90011854009
$ws.Range("S11").Value = 0.09354030474740357
$ws.Range("T11").Value = 0.09354030474740356
$ws.Range("G12").Value = 21.305189
$ws.Range("H12").Value = 63.915567
$ws.Range("I12").Value = 0.1251171236325075
$ws.Range("J12").Value = 0.1251171236325075
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 633.157918413476
$ws.Range("R12").Value = 5698.421265721284
$ws.Range("S12").Value = 0.02991889594284125
$ws.Range("T12").Value = 0.02991889594284125
$ws.Range("G13").Value = 21.305189
$ws.Range("H13").Value = 63.915567
$ws.Range("I13").Value = 0.1251171236325075
$ws.Range("J13").Value = 0.1251171236325075
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 5.878975157848999
$ws.Range("R13").Value = 52.91077642064099
$ws.Range("S13").Value = 0.0002778018577718688
$ws.Range("T13").Value = 0.0002778018577718687
$ws.Range("G14").Value = 0.500358
$ws.Range("H14").Value = 1.501074
$ws.Range("I14").Value = 0.002938408748521978
$ws.Range("J14").Value = 0.002938408748521978
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 0.685928940394
$ws.Range("R14").Value = 6.173360463546
$ws.Range("S14").Value = 0.00003241250878336114
$ws.Range("T14").Value = 0.00003241250878336114
$ws.Range("G15").Value = 0.500358
$ws.Range("H15").Value = 1.501074
$ws.Range("I15").Value = 0.002938408748521978
$ws.Range("J15").Value = 0.002938408748521978
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 46.490125279978
$ws.Range("R15").Value = 418.411127519802
$ws.Range("S15").Value = 0.002196818803287845
$ws.Range("T15").Value = 0.002196818803287845
$ws.Range("G16").Value = 0.500358
$ws.Range("H16").Value = 1.501074
$ws.Range("I16").Value = 0.002938408748521978
$ws.Range("J16").Value = 0.002938408748521978
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 14.869881217272
$ws.Range("R16").Value = 133.828930955448
$ws.Range("S16").Value = 0.0007026531863279018
$ws.Range("T16").Value = 0.0007026531863279017
$ws.Range("G17").Value = 0.500358
$ws.Range("H17").Value = 1.501074
$ws.Range("I17").Value = 0.002938408748521978
$ws.Range("J17").Value = 0.002938408748521978
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 0.138069286878
$ws.Range("R17").Value = 1.242623581902
$ws.Range("S17").Value = 0.000006524250122869913
$ws.Range("T17").Value = 0.000006524250122869912
